$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at positions 17-18, shifting the existing
# rows 17-29 ("lif"/"alif"/"plif" blocks) down to rows 19-31.
$ws.Range("A17:A18").EntireRow.Insert()

# New row 17: a new "tclif" result (no Neuron Model / Network label,
# same pattern as the other blank-labelled rows in that block).
$ws.Range("A17").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "Valve"
$ws.Range("D17").Value = "adam"
$ws.Range("E17").Value = 0.0005
$ws.Range("F17").Value = 256
$ws.Range("G17").Value = 200
$ws.Range("H17").Value = 93.375

# New row 18: another new "tclif" result.
$ws.Range("A18").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "Valve"
$ws.Range("D18").Value = "adam"
$ws.Range("E18").Value = 0.0005
$ws.Range("F18").Value = 256
$ws.Range("G18").Value = 200
$ws.Range("H18").Value = 92.875
